$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '25.958.56'
$ws.Range('E2').Value = '  +0.45%  '

# Row 3
$ws.Range('D3').Value = '1.736.91'
$ws.Range('E3').Value = '  +0.01%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('E5').Value = '  +4.25%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  -0.11%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5045'
$ws.Range('E7').Value = '  -1.43%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2732'
$ws.Range('E8').Value = '  -0.14%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06191'
$ws.Range('E9').Value = '  +1.42%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07263'
$ws.Range('E10').Value = '  +1.29%  '

# Row 11
$ws.Range('D11').Value = '1.737.21'
$ws.Range('E11').Value = '  -0.02%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.6575'
$ws.Range('E12').Value = '  +3.40%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.27'
$ws.Range('E13').Value = '  +2.44%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.751'
$ws.Range('E14').Value = '  +3.54%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.62'
$ws.Range('E15').Value = '  +0.60%  '

# Row 16
$ws.Range('E16').Value = '  -0.03%  '

# Row 17
$ws.Range('E17').Value = '  -0.10%  '

# Row 18
$ws.Range('D18').Value = '25.973.97'
$ws.Range('E18').Value = '  +0.47%  '

# Row 19
$ws.Range('E19').Value = '  +1.74%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006855'
$ws.Range('E20').Value = '  +1.64%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.618'
$ws.Range('E21').Value = '  +8.53%  '

# Row 22
$ws.Range('D22').Value = '1.964.05'
$ws.Range('E22').Value = '  +0.17%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.811'
$ws.Range('E23').Value = '  +1.66%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.453'
$ws.Range('E24').Value = '  +4.42%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '134.36'
$ws.Range('E25').Value = '  -3.30%  '

# Row 26
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '15.28'
$ws.Range('E26').Value = '  +1.11%  '

# Row 27
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.458'
$ws.Range('E27').Value = '  -3.57%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.792'
$ws.Range('E28').Value = '  +2.30%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '105.48'
$ws.Range('E29').Value = '  -0.25%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.990'
$ws.Range('E30').Value = '  -0.41%  '

# Row 31
$ws.Range('E31').Value = '  -2.13%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.715'
$ws.Range('E32').Value = '  +1.93%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04762'
$ws.Range('E33').Value = '  +4.36%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.657'
$ws.Range('E34').Value = '  -0.41%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9975'
$ws.Range('E35').Value = '  +1.34%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6114'
$ws.Range('E36').Value = '  -0.71%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.744'
$ws.Range('E37').Value = '  +2.14%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01614'
$ws.Range('E38').Value = '  +0.37%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8794'
$ws.Range('E39').Value = '  +19.51%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.953'
$ws.Range('E40').Value = '  +1.82%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9994'
$ws.Range('E41').Value = '  -0.13%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '100.85'

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3939'
$ws.Range('E43').Value = '  +2.81%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.022'
$ws.Range('E44').Value = '  +1.47%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1184'
$ws.Range('E45').Value = '  +5.54%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.387'
$ws.Range('E46').Value = '  +3.84%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.77'
$ws.Range('E47').Value = '  +1.82%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05289'
$ws.Range('E48').Value = '  +0.55%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.96'
$ws.Range('E49').Value = '  +1.53%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.3496'
$ws.Range('E50').Value = '  +2.64%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.646'
$ws.Range('E51').Value = '  +1.36%  '
